# Update "想去人数" (want-to-go count) figures in the F column on the
# "展览" and "全部类型" sheets to match the latest scrape (output
# generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$exhibitSheet = $wb.Worksheets.Item("展览")
$exhibitSheet.Range("F3").Value = 131
$exhibitSheet.Range("F4").Value = 198
$exhibitSheet.Range("F5").Value = 3483
$exhibitSheet.Range("F6").Value = 356
$exhibitSheet.Range("F8").Value = 427

$allTypesSheet = $wb.Worksheets.Item("全部类型")
$allTypesSheet.Range("F3").Value = 131
$allTypesSheet.Range("F4").Value = 198
$allTypesSheet.Range("F5").Value = 3483
$allTypesSheet.Range("F6").Value = 356
$allTypesSheet.Range("F10").Value = 427
